$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume/1h (E) updates for rows whose coin identity
# (name/link in columns B/C) does not change. D-column values that look
# numeric are prefixed with "'" so Excel stores them as text, matching
# the original inlineStr (text) cell type used throughout the sheet.
$updates = [ordered]@{
    "D2" = "68.104.53"
    "E2" = "  +7.13%  "
    "D3" = "3.656.23"
    "E3" = "  +5.20%  "
    "E4" = "  -0.12%  "
    "D5" = "'420.48"
    "E5" = "  +1.32%  "
    "D6" = "'129.64"
    "E6" = "  +0.37%  "
    "D7" = "'0.652"
    "E7" = "  +2.74%  "
    "D8" = "3.648.75"
    "E8" = "  +5.19%  "
    "E9" = "  -0.11%  "
    "D10" = "'0.766"
    "E10" = "  +1.82%  "
    "D11" = "'0.196"
    "E11" = "  +26.16%  "
    "D12" = "'0.0000442"
    "E12" = "  +94.77%  "
    "D13" = "'42.00"
    "E13" = "  -0.94%  "
    "D14" = "'9.82"
    "E14" = "  +1.81%  "
    "D15" = "4.218.63"
    "E15" = "  +4.65%  "
    "E16" = "  +0.34%  "
    "D17" = "3.656.19"
    "E17" = "  +5.25%  "
    "D18" = "'20.05"
    "E18" = "  -0.94%  "
    "E19" = "  +2.17%  "
    "D20" = "67.940.28"
    "E20" = "  +6.96%  "
    "E21" = "  +0.52%  "
    "D22" = "'458.37"
    "E22" = "  +0.50%  "
    "E23" = "  -0.83%  "
    "D24" = "'13.48"
    "E24" = "  +2.46%  "
    "D25" = "'3.03"
    "E25" = "  -7.02%  "
    "D26" = "'10.04"
    "E26" = "  -0.94%  "
    "D27" = "'35.73"
    "E27" = "  +6.82%  "
    "E28" = "  -1.80%  "
    "D29" = "'4.98"
    "E29" = "  +4.45%  "
    "D33" = "'7.15"
    "E33" = "  -4.97%  "
    "E34" = "  -6.90%  "
    "D35" = "'40.24"
    "E35" = "  +0.97%  "
    "E36" = "  +0.00%  "
    "E39" = "  +1.57%  "
    "E40" = "  +9.60%  "
    "D41" = "'0.997"
    "E41" = "  -0.30%  "
    "D42" = "'148.56"
    "E42" = "  +1.52%  "
    "E43" = "  -3.30%  "
    "D44" = "'2.93"
    "E44" = "  -5.47%  "
    "D45" = "'2.68"
    "E45" = "  +14.64%  "
    "E46" = "  -2.20%  "
    "D49" = "'0.301"
    "E49" = "  -4.20%  "
    "E50" = "  -2.29%  "
    "D51" = "'2.65"
    "E51" = "  +13.54%  "
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}

# Row re-orderings: coins shift between rows 30-32, 37-38 and 47-48,
# so name (B), link (C), price (D) and volume (E) all change together.
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'12.24"
$ws.Range("E30").Value = "  -1.76%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.119"
$ws.Range("E31").Value = "  +6.09%  "

$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'2.72"
$ws.Range("E32").Value = "  +1.96%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "'56.12"
$ws.Range("E37").Value = "  -2.52%  "

$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0799"
$ws.Range("E38").Value = "  +23.43%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.170"
$ws.Range("E47").Value = "  +22.27%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'4.25"
$ws.Range("E48").Value = "  -5.51%  "

